# TRE - Projeto PUC PR : merge the two trailing runs of the "Obs:" bullet
# on slide 3 into a single run (the stray run break before
# "contratação- manter..." is removed), matching the author's edit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# Locate the textbox shape that holds the "Obs:" bullet point.
$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -like "Obs:*") {
            $target = $sh
            break
        }
    }
}

$tr = $target.TextFrame.TextRange
$fullText = $tr.Text

# The second run starts right after "Obs" (3 characters) and the third run
# starts at the "contratação- manter..." marker. Re-writing the combined
# span (run2 + run3) as a single assignment merges them into one run while
# leaving the leading "Obs" run (and its formatting) untouched.
$marker = "contratação- manter"
$markerPos = $fullText.IndexOf($marker)

if ($markerPos -ge 0) {
    $prefixLen = 3
    $mergedText = $fullText.Substring($prefixLen)
    $mergedLen = $fullText.Length - $prefixLen
    $span = $tr.Characters($prefixLen + 1, $mergedLen)
    $span.Text = $mergedText
}
